$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ROKU")

# Row 4 updates
$ws.Range("B4").Value = 54000000.0
$ws.Range("C4").Value = 63000000.0
$ws.Range("D4").Value = 45000000.0
$ws.Range("E4").Value = 44000000.0
$ws.Range("F4").Value = 50000000.0

# Row 12 updates
$ws.Range("B12").Value = 112000000.0
$ws.Range("C12").Value = 122000000.0
$ws.Range("D12").Value = 131000000.0
$ws.Range("E12").Value = 103000000.0
$ws.Range("F12").Value = 115000000.0
